$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet4 = $wb.Worksheets.Item("全部类型")

$sheet1.Range("F6").Value = 3245
$sheet1.Range("F7").Value = 792
$sheet1.Range("F8").Value = 2008
$sheet1.Range("F9").Value = 1931
$sheet1.Range("F10").Value = 988
$sheet1.Range("F16").Value = 64
$sheet1.Range("F18").Value = 1410
$sheet1.Range("F20").Value = 608
$sheet1.Range("F21").Value = 308
$sheet1.Range("F22").Value = 10565
$sheet1.Range("F23").Value = 9730
$sheet1.Range("F27").Value = 137
$sheet1.Range("F28").Value = 398

$sheet4.Range("F8").Value = 3245
$sheet4.Range("F9").Value = 792
$sheet4.Range("F10").Value = 2008
$sheet4.Range("F11").Value = 1931
$sheet4.Range("F12").Value = 988
$sheet4.Range("F19").Value = 64
$sheet4.Range("F22").Value = 1410
$sheet4.Range("F24").Value = 608
$sheet4.Range("F25").Value = 308
$sheet4.Range("F26").Value = 10565
$sheet4.Range("F27").Value = 9730
$sheet4.Range("F33").Value = 137
$sheet4.Range("F34").Value = 398
